$d = $word.ActiveDocument

# --- Mint the built-in "List" paragraph style (wdStyleList = -48) into the
# document's style sheet without leaving any visible trace on body content.
# Word only serializes a built-in style definition once something forces it
# to resolve concrete pPr/rPr for that style; applying it to a paragraph and
# immediately restoring the paragraph's original style does exactly that.
$mintParagraph = $d.Paragraphs(1).Range
$originalStyle = $mintParagraph.ParagraphStyle
$mintParagraph.Style = $d.Styles(-48)
$mintParagraph.Style = $originalStyle

# Configure the newly minted "List" style so it matches a classic
# "List Paragraph"-style hanging indent with contextual spacing.
$listStyle = $d.Styles("List")
$listStyle.ParagraphFormat.LeftIndent = 18            # 360 twips (0.25")
$listStyle.ParagraphFormat.FirstLineIndent = -18       # -360 twips hanging
$listStyle.NoSpaceBetweenParagraphsOfSameStyle = $true # <w:contextualSpacing/>
$listStyle.UnhideWhenUsed = $true

# --- Re-base "List 1" on the new "List" style instead of "Heading 4".
$list1 = $d.Styles("List1")
$list1.BaseStyle = $listStyle

# Single line spacing (line="240" line-rule="auto") instead of the
# multiple/auto 259 inherited previously.
$list1.ParagraphFormat.LineSpacingRule = 0
$list1.ParagraphFormat.SpaceAfter = 8
$list1.NoSpaceBetweenParagraphsOfSameStyle = $true

# Swap the old "no bold" run-property override for explicit Times New Roman,
# 12 pt (ascii/hAnsi/cs + sz/szCs), matching the new non-heading lineage.
$list1.Font.Name = "Times New Roman"
$list1.Font.NameBi = "Times New Roman"
$list1.Font.Size = 12
$list1.Font.SizeBi = 12
